$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("제출답안")
$ws.Activate()
$ws.Range("B2").Value = "서현"
$ws.Range("B2").Select()
